# Replace the 2018-11-18 (EAST RIDGE) daily accident report data with the
# 2018-11-19 data: rows 2-5 get new values and rows 6-20 are newly added.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row data -------------------------------------------------------------
# row, Accident#, Problem, Latitude, Longitude, Date, Time, Address, City, Hour, Month
$rows = @(
    @(2,  52,  "Injuries",    "35.035569",          "-85.152028",          "2018-11-19", "19:36:24", "2126 Gunbarrel Rd",                 "CHATTANOOGA",     "19", 11),
    @(3,  16,  "Injuries",    "35.084984",          "-85.064733",          "2018-11-19", "18:19:25", "Lee Hwy / Weir Way",                 "CHATTANOOGA",     "18", 11),
    @(4,  17,  "Injuries",    "35.084984",          "-85.064733",          "2018-11-19", "18:19:25", "Lee Hwy / Weir Way",                 "CHATTANOOGA",     "18", 11),
    @(5,  23,  "Injuries",    "35.015214",          "-85.16369899999999",  "2018-11-19", "18:05:06", "Gunbarrel Rd / E Brainerd Rd",       "CHATTANOOGA",     "18", 11),
    @(6,  27,  "Injuries",    "35.007795",          "-85.2311",            "2018-11-19", "17:49:33", "301-329 S Moore Rd",                 "CHATTANOOGA",     "17", 11),
    @(7,  28,  "Injuries",    "35.007795",          "-85.2311",            "2018-11-19", "17:49:33", "301-329 S Moore Rd",                 "CHATTANOOGA",     "17", 11),
    @(8,  29,  "Injuries",    "35.007795",          "-85.2311",            "2018-11-19", "17:49:07", "301-329 S Moore Rd",                 "CHATTANOOGA",     "17", 11),
    @(9,  30,  "Injuries",    "35.007795",          "-85.2311",            "2018-11-19", "17:49:07", "301-329 S Moore Rd",                 "CHATTANOOGA",     "17", 11),
    @(10, 31,  "Injuries",    "35.007795",          "-85.2311",            "2018-11-19", "17:49:06", "301-329 S Moore Rd",                 "CHATTANOOGA",     "17", 11),
    @(11, 40,  "Injuries",    "35.032474",          "-85.263942",          "2018-11-19", "17:12:54", "101 Glenwood Dr",                    "CHATTANOOGA",     "17", 11),
    @(12, 54,  "Injuries",    "35.166365",          "-85.260385",          "2018-11-19", "16:01:29", "6000 Dayton Blvd",                   "CHATTANOOGA",     "16", 11),
    @(13, 55,  "Injuries",    "35.166365",          "-85.260385",          "2018-11-19", "16:00:43", "340-799 MONTLAKE RD",                "CHATTANOOGA",     "16", 11),
    @(14, 98,  "Injuries",    "35.235301",          "-85.22711200000001",  "2018-11-19", "08:09:06", "Bonny Oaks Dr / Hickory Valley Rd",  "HAMILTON COUNTY", "8",  11),
    @(15, 106, "Injuries",    "35.071531",          "-85.156379",          "2018-11-19", "07:24:55", "Wilder St / N Hickory St",           "CHATTANOOGA",     "7",  11),
    @(16, 111, "No Injuries", "35.075148",          "-85.255576",          "2018-11-19", "06:45:03", "3507 Dayton Blvd",                   "CHATTANOOGA",     "6",  11),
    @(17, 114, "Injuries",    "35.113754",          "-85.295055",          "2018-11-19", "05:22:55", "729 E 49th St",                      "RED BANK",        "5",  11),
    @(18, 117, "Injuries",    "34.992135",          "-85.30683399999999",  "2018-11-19", "04:18:23", "5027 Shoals Ln",                     "CHATTANOOGA",     "4",  11),
    @(19, 118, "No Injuries", "34.992135",          "-85.30683399999999",  "2018-11-19", "04:17:00", "5027 Shoals Ln",                     "CHATTANOOGA",     "4",  11),
    @(20, 120, "Injuries",    "35.133212",          "-85.149918",          "2018-11-19", "03:49:18", "6000 Dayton Blvd",                   "HAMILTON COUNTY", "3",  11)
)

# Force the Date (F) and Hour (P) columns to be stored as literal text
# (matching the source data, which keeps these as shared strings rather than
# numeric/date values) instead of letting autodetection coerce them.
$ws.Range("F2:F20").NumberFormat = "@"
$ws.Range("P2:P20").NumberFormat = "@"

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 11).Value = $row[8]
    $ws.Cells.Item($r, 16).Value = $row[9]
    $ws.Cells.Item($r, 25).Value = $row[10]
}

# Put the Date/Hour columns' number format back to General now that the
# text values are committed, so no lingering custom format is left applied.
$ws.Range("F2:F20").Style = "Normal"
$ws.Range("P2:P20").Style = "Normal"

# The Accident# column (A) carries the bold/bordered header-like style on
# every data row; copy that formatting (without touching the new values)
# onto the newly added rows 6-20 so it matches rows 2-5.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A6:A20").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
